# Slide 1, title placeholder ("Rectangle 62" / ctrTitle) currently reads
# "kathara lab" as a single run. The commit retitles it to
# "Lab Random Loadbalancer", split across two runs:
#   1) "Lab Random "
#   2) "Loadbalancer"
# (the second run is flagged by PowerPoint's proofing tools as a
# possible spelling error, consistent with it being a compound/product
# word not in the dictionary).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Replace the whole title text with the first new run...
$tr.Text = "Lab Random "

# ...then append the second run right after it, so the paragraph ends
# up with two separate runs (matching the authored split) instead of
# one merged run.
$tr2 = $tr.InsertAfter("Loadbalancer")
